# Update "Free Space" figures on the Server Spaces sheet to reflect the
# latest capacity scan. Values are stored as plain text (e.g. "443.85 GB"),
# so we assign them as strings to keep the existing text formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value   = "443.85 GB"
$ws.Range("D11").Value  = "204.94 GB"
$ws.Range("D18").Value  = "30.33 GB"
$ws.Range("D117").Value = "30.33 GB"
$ws.Range("D28").Value  = "75.09 GB"
$ws.Range("D127").Value = "75.09 GB"
$ws.Range("D41").Value  = "97.21 GB"
$ws.Range("D43").Value  = "8.13 GB"
$ws.Range("D44").Value  = "1021.72 GB"
$ws.Range("D45").Value  = "527.66 GB"
$ws.Range("D48").Value  = "69.63 GB"
$ws.Range("D67").Value  = "98.8 GB"
$ws.Range("D70").Value  = "300.98 GB"
$ws.Range("D71").Value  = "1393.48 GB"
$ws.Range("D74").Value  = "899.28 GB"
$ws.Range("D78").Value  = "209.11 GB"
$ws.Range("D143").Value = "209.11 GB"
$ws.Range("D101").Value = "11.93 GB"
$ws.Range("D104").Value = "18.15 GB"
$ws.Range("D105").Value = "906.54 GB"
$ws.Range("D113").Value = "9.74 GB"
$ws.Range("D158").Value = "630.7 GB"
$ws.Range("D160").Value = "153.35 GB"
$ws.Range("D163").Value = "430.92 GB"
$ws.Range("D165").Value = "424.88 GB"
$ws.Range("D169").Value = "2.28 GB"
$ws.Range("D174").Value = "3.9 GB"
$ws.Range("D176").Value = "84.01 GB"
$ws.Range("D178").Value = "56.92 GB"
$ws.Range("D183").Value = "191.64 GB"
$ws.Range("D185").Value = "46.62 GB"
$ws.Range("D191").Value = "163.58 GB"
$ws.Range("D199").Value = "68.37 GB"
$ws.Range("D207").Value = "74.5 GB"
$ws.Range("D223").Value = "124.97 GB"
